$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "328.18"
Set-TextValue "E2" "-1.36%"
Set-TextValue "D3" "44.83"
Set-TextValue "E3" "-2.29%"
Set-TextValue "D4" "5.286"
Set-TextValue "E4" "-4.78%"
Set-TextValue "D5" "0.08372"
Set-TextValue "E5" "0.33%"
Set-TextValue "D6" "1.950"
Set-TextValue "E6" "-4.54%"
Set-TextValue "D7" "0.9714"
Set-TextValue "E7" "-0.91%"
Set-TextValue "D9" "0.1108"
Set-TextValue "E9" "-2.62%"
Set-TextValue "D10" "0.1910"
Set-TextValue "E10" "-1.23%"
Set-TextValue "D11" "0.09702"
Set-TextValue "E11" "-3.49%"
Set-TextValue "D12" "0.04601"
Set-TextValue "E12" "-0.56%"
Set-TextValue "D13" "0.1060"
Set-TextValue "E13" "-0.10%"
Set-TextValue "D14" "0.001290"
Set-TextValue "E14" "2.04%"
Set-TextValue "D15" "0.005878"
Set-TextValue "E15" "-0.40%"
Set-TextValue "E16" "0.03%"
Set-TextValue "D17" "4.423"
Set-TextValue "E17" "-0.19%"
Set-TextValue "D18" "0.3355"
Set-TextValue "E18" "0.23%"
Set-TextValue "D19" "8.351"
Set-TextValue "E19" "-18.95%"
Set-TextValue "E20" "-2.38%"
Set-TextValue "D21" "0.2718"
Set-TextValue "E21" "9.26%"
Set-TextValue "E22" "1.74%"
Set-TextValue "E23" "-4.81%"
Set-TextValue "D24" "0.004447"
Set-TextValue "E24" "0.49%"
Set-TextValue "D25" "0.0001301"
Set-TextValue "E25" "1.63%"
Set-TextValue "D26" "0.0002978"
Set-TextValue "E26" "-20.31%"
Set-TextValue "D38" "0.02717"
Set-TextValue "E38" "-3.98%"
Set-TextValue "D39" "0.05637"
Set-TextValue "E39" "-2.43%"
Set-TextValue "D40" "0.007763"
Set-TextValue "E40" "1.27%"
Set-TextValue "E41" "-1.18%"
Set-TextValue "D42" "0.007319"
Set-TextValue "E42" "-3.09%"
Set-TextValue "D43" "0.002114"
Set-TextValue "E43" "7.11%"
Set-TextValue "D44" "0.007882"
Set-TextValue "E44" "-1.86%"
Set-TextValue "D45" "0.3507"
Set-TextValue "D46" "0.00006958"
Set-TextValue "E46" "-2.86%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "-0.02%"
Set-TextValue "D48" "0.003488"
Set-TextValue "E48" "-0.40%"
Set-TextValue "D49" "0.003529"
Set-TextValue "E49" "39.82%"
Set-TextValue "D50" "0.00002099"
Set-TextValue "E50" "-0.02%"
Set-TextValue "D51" "0.0001999"
Set-TextValue "E51" "-0.02%"
